$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 54412
$ws.Cells.Item(2, 2).Value = "Alana Duarte"
$ws.Cells.Item(2, 3).Value = "Atendimento ao Cliente"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 45087
$ws.Cells.Item(2, 7).Value = 7984.51

# Row 3
$ws.Cells.Item(3, 1).Value = 92464
$ws.Cells.Item(3, 2).Value = "Camila Moreira"
$ws.Cells.Item(3, 3).Value = "Engenharia"
$ws.Cells.Item(3, 4).Value = "Outros"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 45080
$ws.Cells.Item(3, 7).Value = 3206.79

# Row 4
$ws.Cells.Item(4, 1).Value = 28106
$ws.Cells.Item(4, 2).Value = "Davi Lucas Carvalho"
$ws.Cells.Item(4, 3).Value = "Vendas"
$ws.Cells.Item(4, 4).Value = "Doença"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 45097
$ws.Cells.Item(4, 7).Value = 4345.63

# Row 5
$ws.Cells.Item(5, 1).Value = 38701
$ws.Cells.Item(5, 2).Value = "Catarina Jesus"
$ws.Cells.Item(5, 3).Value = "Engenharia"
$ws.Cells.Item(5, 4).Value = "Doença"
$ws.Cells.Item(5, 5).Value = 4
$ws.Cells.Item(5, 6).Value = 45079
$ws.Cells.Item(5, 7).Value = 11845.44

# Row 6
$ws.Cells.Item(6, 1).Value = 81682
$ws.Cells.Item(6, 2).Value = "Gabrielly Souza"
$ws.Cells.Item(6, 3).Value = "TI"
$ws.Cells.Item(6, 4).Value = "Outros"
$ws.Cells.Item(6, 5).Value = 7
$ws.Cells.Item(6, 6).Value = 45102
$ws.Cells.Item(6, 7).Value = 7355.6

# Row 7
$ws.Cells.Item(7, 1).Value = 16630
$ws.Cells.Item(7, 2).Value = "Daniela Santos"
$ws.Cells.Item(7, 3).Value = "TI"
$ws.Cells.Item(7, 4).Value = "Consulta médica"
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = 45094
$ws.Cells.Item(7, 7).Value = 8911.78

# Row 8
$ws.Cells.Item(8, 1).Value = 87459
$ws.Cells.Item(8, 2).Value = "Valentina Ribeiro"
$ws.Cells.Item(8, 3).Value = "Marketing"
$ws.Cells.Item(8, 4).Value = "Doença"
$ws.Cells.Item(8, 5).Value = 5
$ws.Cells.Item(8, 6).Value = 45104
$ws.Cells.Item(8, 7).Value = 11025.51

# Row 9
$ws.Cells.Item(9, 1).Value = 45912
$ws.Cells.Item(9, 2).Value = "Dr. Bruno da Mota"
$ws.Cells.Item(9, 4).Value = "Doença"
$ws.Cells.Item(9, 5).Value = 7
$ws.Cells.Item(9, 6).Value = 45092
$ws.Cells.Item(9, 7).Value = 8250.62

# Row 10
$ws.Cells.Item(10, 1).Value = 21670
$ws.Cells.Item(10, 2).Value = "Mirella da Costa"
$ws.Cells.Item(10, 3).Value = "TI"
$ws.Cells.Item(10, 4).Value = "Consulta médica"
$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(10, 6).Value = 45083
$ws.Cells.Item(10, 7).Value = 11302.07

# Row 11
$ws.Cells.Item(11, 1).Value = 76316
$ws.Cells.Item(11, 2).Value = "Noah da Costa"
$ws.Cells.Item(11, 3).Value = "Financeiro"
$ws.Cells.Item(11, 5).Value = 7
$ws.Cells.Item(11, 6).Value = 45095
$ws.Cells.Item(11, 7).Value = 5598.35
